$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved (these columns store numeric-looking
# and percentage-looking values as literal text, not real numbers, so we
# must force Text format before assigning or Excel will coerce them to
# floating point numbers / dates and lose exact formatting and precision).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '34.742.60'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '1.807.21'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '225.11'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").Value = '0.556'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '32.55'
$ws.Range("E8").Value = '  +4.69%  '
$ws.Range("D9").Value = '0.289'
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("D10").Value = '0.0715'
$ws.Range("E10").Value = '  +8.20%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '2.066.31'
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '1.805.19'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '0.644'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '34.746.49'
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").Value = '  +3.03%  '
$ws.Range("D18").Value = '69.80'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '255.18'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '0.0₃0810'
$ws.Range("E20").Value = '  +9.03%  '
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '10.88'
$ws.Range("E22").Value = '  +4.57%  '
$ws.Range("D23").Value = '4.27'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = '161.32'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").Value = '16.53'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '7.19'
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +3.87%  '
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").Value = '1.21'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("E34").Value = '  +3.45%  '
$ws.Range("D35").Value = '1.449.41'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +3.80%  '
$ws.Range("D38").Value = '0.641'
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("D39").Value = '85.60'
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("D40").Value = '0.962'
$ws.Range("E40").Value = '  +6.69%  '
$ws.Range("D41").Value = '2.80'
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("D42").Value = '2.34'
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("E44").Value = '  +6.83%  '
$ws.Range("B45").Value = 'Swop.fi'
$ws.Range("C45").Value = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop'
$ws.Range("D45").Value = '321.94'
$ws.Range("E45").Value = '  +510.42%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.06'
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").Value = '0.0493'
$ws.Range("E47").Value = '  -4.15%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.962.09'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '106.09'
$ws.Range("E49").Value = '  +8.97%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '12.14'
$ws.Range("E50").Value = '  +3.63%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.05%  '
